# Add new "Block" property rows (Grass1-5, Crack1-5, Treasure1 , Treasure2)
# to the Property sheet, following the existing row pattern (rows 3-10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    "Grass1",
    "Grass2",
    "Grass3",
    "Grass4",
    "Grass5",
    "Crack1",
    "Crack2",
    "Crack3",
    "Crack4",
    "Crack5",
    "Treasure1 ",
    "Treasure2"
)

$startRow = 11
for ($i = 0; $i -lt $names.Count; $i++) {
    $row = $startRow + $i

    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = "string"
    $ws.Cells.Item($row, 3).Value = $false
    $ws.Cells.Item($row, 4).Value = $false
    $ws.Cells.Item($row, 5).Value = $false
    $ws.Cells.Item($row, 6).Value = $true
    $ws.Cells.Item($row, 7).Value = 0
    $ws.Cells.Item($row, 8).Value = 0
    $ws.Cells.Item($row, 9).Value = "Friend"

    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 9).NumberFormat = "@"
    if ($row -le 20) {
        $ws.Cells.Item($row, 1).NumberFormat = "@"
    }
}

$ws.Range("H24").Select()
